# Applies the "Update functional requirements (belongs to DatNT)" edit.
$d = $word.ActiveDocument

# Helper: isolate [start,end) into its own run so that newly inserted /
# adjacent text does not get silently re-coalesced into the neighbouring
# run. Toggling a character-format property on and back off over the
# exact span forces a fresh run boundary without altering the visible
# formatting.
function Isolate-Range($doc, $s, $e) {
    $r = $doc.Range($s, $e)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark next to the smart-watch sentence.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) "Approve bus route and bus time change from background handler." ->
#    "...handler and write to official database." with a new "_GoBack"
#    bookmark placed right before the final period.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Approve bus route and bus time change from background handler.")
$periodPos = $rng.End - 1

$insertRng = $d.Range($periodPos, $periodPos)
$insertRng.InsertBefore(" and write to official database")

$newTextLen = (" and write to official database").Length
$newTextStart = $periodPos
$newTextEnd = $periodPos + $newTextLen

# Isolate the newly inserted text into its own run (separates it both from
# "...handler" on the left and the trailing "." on the right).
Isolate-Range $d $newTextStart $newTextEnd

$gbRng = $d.Range($newTextEnd, $newTextEnd)
$d.Bookmarks.Add("_GoBack", $gbRng)

# ---------------------------------------------------------------------
# 3) "Periodically detect change from official bus website." ->
#    "...website and write to temporary database."
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Periodically detect change from official bus website.")
$periodPos2 = $rng2.End - 1

$insertRng2 = $d.Range($periodPos2, $periodPos2)
$insertRng2.InsertBefore(" and write to temporary database")

$newTextLen2 = (" and write to temporary database").Length
$newTextStart2 = $periodPos2
$newTextEnd2 = $periodPos2 + $newTextLen2

Isolate-Range $d $newTextStart2 $newTextEnd2

# ---------------------------------------------------------------------
# 4) Insert a new bullet "Synchronize data from server to mobile." right
#    after the "Mobile Component:" paragraph (before the first "Find
#    optimize path..." bullet).
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Find optimize path between two point when using bus.")
$origStart3 = $rng3.Start
$p3 = $rng3.Paragraphs(1)
$p3.Range.InsertParagraphBefore()

$insertPointRng3 = $d.Range($origStart3, $origStart3)
$insertPointRng3.InsertAfter("Synchronize data from server to mobile.")

Write-Host "Stage 1-4 complete"
